$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 151; existing rows 151-158 shift down to 152-159.
$ws.Rows(151).Insert()

$ws.Range("A151").Value = 5
$ws.Range("B151").Value = "Macroferia Regional de Talca"
$ws.Range("C151").Value = "Maule"
$ws.Range("D151").Value = 44931
$ws.Range("D151").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E151").Value = 7
$ws.Range("F151").Value = "Fruta"
$ws.Range("G151").Value = 100108
$ws.Range("H151").Value = "Tropicales y subtropicales"
$ws.Range("I151").Value = 100108002
$ws.Range("J151").Value = "Mango"
$ws.Range("K151").Value = "Sin especificar"
$ws.Range("L151").Value = "Primera"
$ws.Range("M151").Value = 248
$ws.Range("N151").Value = 6000
$ws.Range("O151").Value = 6000
$ws.Range("P151").Value = 6000
$ws.Range("Q151").Value = "`$/bandeja 4 kilos"
$ws.Range("R151").Value = "Perú"
$ws.Range("S151").Value = 1500
$ws.Range("T151").Value = 4
